$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F values (JAGS LLL model results) that were missing from the upload
$ws.Range("F1").Value = "Jags LLL"

$ws.Range("F3").Value = "5.6 (4.4,7.2)"
$ws.Range("F4").Value = "1.44 (1.22, 1.9)"
$ws.Range("F5").Value = "3.2 (1.7, 4.5)"
$ws.Range("F6").Value = "10.5(7.3, 16.9)"

$ws.Range("F10").Value = "3.9 (0.1,7.7)"
$ws.Range("F11").Value = "8.5 (3.8,12)"
$ws.Range("F12").Value = "33.4(12,592)"

$ws.Range("F16").Value = "6.9 (3.6, 8.8)"
$ws.Range("F17").Value = "9.0 (7.4, 10.1)"
$ws.Range("F18").Value = "12.4 (8.9,19.5)"

# Make the new header cell bold to match the other header cells in row 1
$ws.Range("F1").Font.Bold = $true

# Update the active selection to match the edited workbook
$ws.Range("E10").Select()
